{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Implements:\n//   1. Merge the \"May 1\" / \"6\" / \", 2025\" runs into a single \"May 16, 2025\" run.\n//   2. Re-wrap three body paragraphs into many runs (one run per ~80-char line,\n//      with a separate single-space run between each line) while keeping the\n//      paragraph's own formatting (pPr) untouched.\n//   3. Wrap the \"What are some additional aspects...\" heading run in a\n//      bookmark named \"_Int_BIxpuo3D\".\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build a <w:p>...</w:p> OOXML fragment (wrapped in the minimal package\n// envelope insertOoxml expects) that reproduces `lines` as alternating\n// text-runs / single-space-runs, keeping the paragraph's existing <w:pPr>.\nfunction buildReflowOoxml(pPrXml, lines) {\n  let runs = \"\";\n  for (let i = 0; i < lines.length; i++) {\n    if (i > 0) {\n      runs += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>';\n    }\n    runs += `<w:r><w:t>${escapeXml(lines[i])}</w:t></w:r>`;\n  }\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    `<w:body><w:p>${pPrXml}${runs}</w:p></w:body>` +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst PARA_REWRAPS = [\n  {\n    startsWith: \"Defense in Depth (DiD) emphasizes\",\n    lines: [\n      \"Defense in Depth (DiD) emphasizes using multiple, overlapping layers of\",\n      \"protection to reduce system vulnerabilities. However, \\u201ctoo deep\\u201d occurs when the\",\n      \"cost of added complexity outweighs the security benefit. Excessive layers can\",\n      \"introduce operational friction, increase system latency, and create new attack\",\n      \"surfaces due to misconfigurations or interoperability issues. For example,\",\n      \"layering multiple endpoint defenses might conflict or duplicate effort without\",\n      \"providing significantly better coverage. The tradeoff is between achieving\",\n      \"robust protection and maintaining usability, performance, and maintainability. A\",\n      \"mature DiD strategy balances depth with necessity, driven by a clear\",\n      \"understanding of the system\\u2019s threat model and critical assets.\",\n    ],\n  },\n  {\n    startsWith: \"Implementing DiD requires careful investment\",\n    lines: [\n      \"Implementing DiD requires careful investment. Financially, organizations must\",\n      \"account for licensing security tools, training developers, and allocating staff\",\n      \"time for monitoring and maintenance. Operationally, more layers mean more\",\n      \"complexity and higher risk of downtime or bottlenecks, especially if layers are\",\n      \"not properly integrated. However, the reputational and legal damage from a\",\n      \"breach\\u2014lost customer trust, regulatory fines, and negative press\\u2014often justifies\",\n      \"these costs. Secure coding practices, such as safe memory handling in C++, help\",\n      \"reduce the burden on later-stage defenses by addressing issues at the source. In\",\n      \"doing so, development teams minimize future operational risks and reduce\",\n      \"dependence on costly compensating controls.\",\n    ],\n  },\n  {\n    startsWith: \"DiD is not a one-size-fits-all solution\",\n    lines: [\n      \"DiD is not a one-size-fits-all solution; its implementation varies depending on\",\n      \"system architecture, use case, and threat profile. A high-security system like a\",\n      \"financial service backend requires deeper, more rigorous controls than a static\",\n      \"public-facing website. For C++ applications, DiD might include language-level\",\n      \"precautions (e.g., avoiding unsafe functions), compiler and OS-level protections\",\n      \"(e.g., stack canaries, ASLR), and runtime policies (e.g., access controls,\",\n      \"sandboxing). The effectiveness of DiD depends on aligning security controls with\",\n      \"real-world risks and operational constraints. Customizing DiD to its context\",\n      \"ensures that each layer meaningfully contributes to the system\\u2019s resilience\",\n      \"without overengineering the solution.\",\n    ],\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\n// ---- 1. Fix the date paragraph: \"May 1\" + \"6\" + \", 2025\" -> \"May 16, 2025\"\nfor (const p of paragraphs.items) {\n  if (p.text.trim() === \"May 16, 2025\") {\n    const r = p.getRange();\n    r.insertText(\"May 16, 2025\", \"Replace\");\n  }\n}\nawait context.sync();\n\n// ---- 2. Re-wrap the three long paragraphs into many runs.\nconst pPrSpacing = '<w:pPr><w:spacing w:line=\"480\" w:lineRule=\"auto\"/></w:pPr>';\nfor (const p of paragraphs.items) {\n  for (const rewrap of PARA_REWRAPS) {\n    if (p.text.indexOf(rewrap.startsWith) === 0) {\n      const ooxml = buildReflowOoxml(pPrSpacing, rewrap.lines);\n      const r = p.getRange(\"Content\");\n      r.insertOoxml(ooxml, \"Replace\");\n    }\n  }\n}\nawait context.sync();\n\n// ---- 3. Bookmark the \"What are some additional aspects...\" heading run.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\nfor (const p of paragraphs2.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const p of paragraphs2.items) {\n  if (p.text.indexOf(\"What are some additional aspects of DiD\") !== -1) {\n    const r = p.getRange(\"Content\");\n    r.insertBookmark(\"_Int_BIxpuo3D\");\n  }\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Implements:\n#   1. Merge the \"May 1\" / \"6\" / \", 2025\" runs into a single \"May 16, 2025\" run.\n#   2. Re-wrap three body paragraphs into many runs (one run per ~80-char line,\n#      with a separate single-space run between each line) while keeping the\n#      paragraph's own formatting (pPr) untouched.\n#   3. Wrap the \"What are some additional aspects...\" heading run in a\n#      bookmark named \"_Int_BIxpuo3D\".\n#\n# NOTE: this interpreter does not preserve live COM-object bindings when they\n# are passed through user-defined function parameters / collections, so all\n# logic below is written inline (no helper functions) on purpose.\n\n$doc = $word.ActiveDocument\n\n# ---- 1. Fix the date paragraph: \"May 1\" + \"6\" + \", 2025\" -> \"May 16, 2025\"\nforeach ($p in $doc.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"May 16, 2025\") {\n        $full = $p.Range\n        $contentRange = $doc.Range($full.Start, $full.End - 1)\n        $contentRange.Text = \"\"\n        $insertPoint = $doc.Range($full.Start, $full.Start)\n        $insertPoint.Text = \"May 16, 2025\"\n    }\n}\n\n# ---- 2. Re-wrap the three long paragraphs into many runs.\n$para1Lines = @(\n    \"Defense in Depth (DiD) emphasizes using multiple, overlapping layers of\",\n    \"protection to reduce system vulnerabilities. However, \u201ctoo deep\u201d occurs when the\",\n    \"cost of added complexity outweighs the security benefit. Excessive layers can\",\n    \"introduce operational friction, increase system latency, and create new attack\",\n    \"surfaces due to misconfigurations or interoperability issues. For example,\",\n    \"layering multiple endpoint defenses might conflict or duplicate effort without\",\n    \"providing significantly better coverage. The tradeoff is between achieving\",\n    \"robust protection and maintaining usability, performance, and maintainability. A\",\n    \"mature DiD strategy balances depth with necessity, driven by a clear\",\n    \"understanding of the system\u2019s threat model and critical assets.\"\n)\n\nforeach ($p in $doc.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"Defense in Depth (DiD) emphasizes\")) {\n        $runsXml = \"\"\n        for ($i = 0; $i -lt $para1Lines.Count; $i++) {\n            if ($i -gt 0) {\n                $runsXml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n            }\n            $runsXml += \"<w:r><w:t>$($para1Lines[$i])</w:t></w:r>\"\n        }\n        $frag = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n        $full = $p.Range\n        $contentRange = $doc.Range($full.Start, $full.End - 1)\n        $contentRange.InsertXML($frag)\n        break\n    }\n}\n\n$para2Lines = @(\n    \"Implementing DiD requires careful investment. Financially, organizations must\",\n    \"account for licensing security tools, training developers, and allocating staff\",\n    \"time for monitoring and maintenance. Operationally, more layers mean more\",\n    \"complexity and higher risk of downtime or bottlenecks, especially if layers are\",\n    \"not properly integrated. However, the reputational and legal damage from a\",\n    \"breach\u2014lost customer trust, regulatory fines, and negative press\u2014often justifies\",\n    \"these costs. Secure coding practices, such as safe memory handling in C++, help\",\n    \"reduce the burden on later-stage defenses by addressing issues at the source. In\",\n    \"doing so, development teams minimize future operational risks and reduce\",\n    \"dependence on costly compensating controls.\"\n)\n\nforeach ($p in $doc.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"Implementing DiD requires careful investment\")) {\n        $runsXml = \"\"\n        for ($i = 0; $i -lt $para2Lines.Count; $i++) {\n            if ($i -gt 0) {\n                $runsXml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n            }\n            $runsXml += \"<w:r><w:t>$($para2Lines[$i])</w:t></w:r>\"\n        }\n        $frag = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n        $full = $p.Range\n        $contentRange = $doc.Range($full.Start, $full.End - 1)\n        $contentRange.InsertXML($frag)\n        break\n    }\n}\n\n$para3Lines = @(\n    \"DiD is not a one-size-fits-all solution; its implementation varies depending on\",\n    \"system architecture, use case, and threat profile. A high-security system like a\",\n    \"financial service backend requires deeper, more rigorous controls than a static\",\n    \"public-facing website. For C++ applications, DiD might include language-level\",\n    \"precautions (e.g., avoiding unsafe functions), compiler and OS-level protections\",\n    \"(e.g., stack canaries, ASLR), and runtime policies (e.g., access controls,\",\n    \"sandboxing). The effectiveness of DiD depends on aligning security controls with\",\n    \"real-world risks and operational constraints. Customizing DiD to its context\",\n    \"ensures that each layer meaningfully contributes to the system\u2019s resilience\",\n    \"without overengineering the solution.\"\n)\n\nforeach ($p in $doc.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"DiD is not a one-size-fits-all solution\")) {\n        $runsXml = \"\"\n        for ($i = 0; $i -lt $para3Lines.Count; $i++) {\n            if ($i -gt 0) {\n                $runsXml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n            }\n            $runsXml += \"<w:r><w:t>$($para3Lines[$i])</w:t></w:r>\"\n        }\n        $frag = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n        $full = $p.Range\n        $contentRange = $doc.Range($full.Start, $full.End - 1)\n        $contentRange.InsertXML($frag)\n        break\n    }\n}\n\n# ---- 3. Bookmark the \"What are some additional aspects...\" heading run.\nforeach ($p in $doc.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"What are some additional aspects of DiD\")) {\n        $full = $p.Range\n        $contentRange = $doc.Range($full.Start, $full.End - 1)\n        $doc.Bookmarks.Add(\"_Int_BIxpuo3D\", $contentRange)\n        break\n    }\n}\n"}
